$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Enum" column (D) to host "Date"
$ws.Columns("D:D").Insert()

# Header for the new "Date" column
$ws.Range("D1").Value = "Date"
$ws.Range("D2").Value = "[No Date]"

# Row3 values: C3 keeps the full DateTime (date + time-of-day),
# D3 gets the date-only value.
$ws.Range("D3").Value = 43831
$ws.Range("D3").NumberFormat = "yyyy-MM-dd"

$ws.Range("C3").Value = 43831.4271412037
$ws.Range("C3").NumberFormat = "yyyy-MM-dd HH:mm:ss"

# Column widths (values compensate for the host's whole-pixel rounding of
# ColumnWidth so the saved character width lands as close as possible to
# the widths used by the original authoring tool: 20.567768 / 12.424911)
$ws.Columns("C:C").ColumnWidth = 19.734435
$ws.Columns("D:D").ColumnWidth = 11.591578

# Re-apply the autofilter over the new full range
$ws.AutoFilterMode = $false
$ws.Range("A1:F3").AutoFilter() | Out-Null

# Fix up the hidden _FilterDatabase defined name so it spans the new range
$wb.Names.Item(1).RefersTo = "=Sheet1!`$A`$1:`$F`$3"
